# "añadido logica para la temporada 2006"
# 1) Extend the existing "Hoja1" table with the 2005 and 2006 seasons.
# 2) Duplicate the whole (now 5-row) table into a new "Hoja2" worksheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# --- Hoja1: append rows for the 2005 and 2006 seasons -----------------
$ws1.Range("A4").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2005.xlsx"
$ws1.Range("B4").Value = "2005"
$ws1.Range("A5").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2006.xlsx"
$ws1.Range("B5").Value = "2006"

# --- Hoja2: new worksheet placed right after Hoja1 ---------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

# file_name column (B) must stay text, same as on Hoja1, so the
# year values ("2003", "2004", ...) are not auto-converted to numbers.
$ws2.Range("B1:B5").NumberFormat = "@"

$ws2.Range("A1").Value = "file_path"
$ws2.Range("B1").Value = "file_name"
$ws2.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2003.xlsx"
$ws2.Range("B2").Value = "2003"
$ws2.Range("A3").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2004.xlsx"
$ws2.Range("B3").Value = "2004"
$ws2.Range("A4").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2005.xlsx"
$ws2.Range("B4").Value = "2005"
$ws2.Range("A5").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2006.xlsx"
$ws2.Range("B5").Value = "2006"

# Match column A's "best fit" width from Hoja1 (~51 characters wide).
$ws2.Columns.Item(1).ColumnWidth = 50.16666667

# Leave Hoja2 with A2:B5 selected (first data row active) ...
[void]$ws2.Range("A2:B5").Select()

# ... but make Hoja1 the active sheet again, selection on B9.
$ws1.Activate()
[void]$ws1.Range("B9").Select()
